# Spéc détaillée formulaire réservation - edits per commit
# "modif spé détaillée formulaire réservation"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: locate a paragraph (by index in $d.Paragraphs) whose text
# matches a given substring, returning its Range.
# ---------------------------------------------------------------------
function Get-ParaRange($matchText) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -match [regex]::Escape($matchText)) {
            return $p.Range
        }
    }
    return $null
}

# =======================================================================
# 1) Untracked fix: merge "comb" + (bookmark _GoBack) + "obox" into a
#    single "combobox" run, dropping the old _GoBack bookmark.
#    NB: the "Chambres" paragraph mentions "combobox" twice ("qui seront
#    dans une « combobox »" and "qui sera aussi une « combobox »") - the
#    broken comb/obox split (with the stray _GoBack bookmark) is the
#    *second* one, so anchor the search right after "sera aussi une".
# =======================================================================
$comboPara = Get-ParaRange("sera aussi une")
$comboAnchor = $comboPara.Duplicate
$comboAnchor.Find.ClearFormatting()
$comboAnchor.Find.Execute("sera aussi une", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($comboAnchor.Find.Found) {
    $comboRng = $d.Range($comboAnchor.End, $comboPara.End)
    $comboRng.Find.ClearFormatting()
    $comboRng.Find.Execute("combobox", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if ($comboRng.Find.Found) {
        # Force a real mutation (identical-text assignment is a no-op),
        # then set the final text so the two split runs merge into one
        # and the bookmark embedded between them is dropped.
        $comboRng.Text = "combobox_TMPFIX"
        $comboRng.Text = "combobox"
    }
}

# =======================================================================
# From here on, edits are tracked changes (Word's "Suivi des modifications").
# =======================================================================
$word.UserName = "59011-14-07"
$d.TrackRevisions = $true

# -----------------------------------------------------------------------
# 2) Paragraph "Les informations d'un client peuvent être retrouvées..."
#    - insert "," after "numéro d'entreprise" (before " ou le numéro de")
#    - insert "," after "bel et bien retrouvées" (before " les différents
#      champs")
# -----------------------------------------------------------------------
$p1 = Get-ParaRange("Les informations d’un client peuvent être retrouvées")
$r1 = $p1.Duplicate
$r1.Find.ClearFormatting()
$r1.Find.Execute("numéro d’entreprise ou le numéro de", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r1.Find.Found) {
    # Collapse to just after "numéro d’entreprise"
    $r1.End = $r1.Start + ("numéro d’entreprise").Length
    $r1.Start = $r1.End
    $r1.InsertAfter(",")
}

$p1b = Get-ParaRange("Les informations d’un client peuvent être retrouvées")
$r2 = $p1b.Duplicate
$r2.Find.ClearFormatting()
$r2.Find.Execute("bel et bien retrouvées les différents champs", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r2.Find.Found) {
    $r2.End = $r2.Start + ("bel et bien retrouvées").Length
    $r2.Start = $r2.End
    $r2.InsertAfter(",")
}

# -----------------------------------------------------------------------
# 3) "Une infobulle..." paragraph: delete "courte " (tracked deletion)
# -----------------------------------------------------------------------
$p2 = Get-ParaRange("infobulle indiquant le rang de fidélité")
$r3 = $p2.Duplicate
$r3.Find.ClearFormatting()
$r3.Find.Execute("courte ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r3.Find.Found) {
    $r3.Delete()
}

# -----------------------------------------------------------------------
# 4) "Chambres" paragraph: bold the "+" in "« + »" (tracked format change)
# -----------------------------------------------------------------------
$p3 = Get-ParaRange("sera aussi une")
$r4 = $p3.Duplicate
$r4.Find.ClearFormatting()
$r4.Find.Execute("+", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r4.Find.Found) {
    $r4.Font.Bold = 1
}

# -----------------------------------------------------------------------
# 5) "Prestations" paragraph: bold the "+" in "« + »" (tracked format
#    change)
# -----------------------------------------------------------------------
$p4 = Get-ParaRange("le réceptionniste pourra ajouter des")
$r5 = $p4.Duplicate
$r5.Find.ClearFormatting()
$r5.Find.Execute("+", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r5.Find.Found) {
    $r5.Font.Bold = 1
}

# -----------------------------------------------------------------------
# 6) "Pour finir..." paragraph:
#    - bold "Valider"
#    - bold "Réinitialiser"
#    - insert the "(pour la sécurité une « Message box » permettant de
#      valider la réinitialisation)," parenthetical + comma after the
#      "Réinitialiser »" closing guillemet
#    - bold "Annuler"
#    - re-add the _GoBack bookmark at the very end of the paragraph
# -----------------------------------------------------------------------
$p5 = Get-ParaRange("Pour finir le réceptionniste pourra soit valider")

$r6 = $p5.Duplicate
$r6.Find.ClearFormatting()
$r6.Find.Execute("Valider", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r6.Find.Found) {
    $r6.Font.Bold = 1
}

$p5b = Get-ParaRange("Pour finir le réceptionniste pourra soit valider")
$r7 = $p5b.Duplicate
$r7.Find.ClearFormatting()
$r7.Find.Execute("Réinitialiser", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r7.Find.Found) {
    $r7.Font.Bold = 1
}

# Insert the parenthetical right after "Réinitialiser »" (i.e. right
# before " ou pour finir annuler...")
$p5c = Get-ParaRange("Pour finir le réceptionniste pourra soit valider")
$r8 = $p5c.Duplicate
$r8.Find.ClearFormatting()
$r8.Find.Execute(" ou pour finir annuler la saisie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r8.Find.Found) {
    $r8.Collapse(1)  # wdCollapseStart

    # " (pour la sécurité une « "
    $r8.InsertAfter(" (pour la sécurité une « ")
    $r8.Collapse(0)  # wdCollapseEnd, move past just-inserted text

    # "Message box" in bold
    $r8.InsertAfter("Message box")
    $r8.Font.Bold = 1
    $r8.Collapse(0)

    # " » permettant de valider la réinitialisation),"
    $r8.InsertAfter(" » permettant de valider la réinitialisation),")
}

$p5d = Get-ParaRange("Pour finir le réceptionniste pourra soit valider")
$r9 = $p5d.Duplicate
$r9.Find.ClearFormatting()
$r9.Find.Execute("Annuler", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($r9.Find.Found) {
    $r9.Font.Bold = 1
}

# Re-add the _GoBack bookmark at the very end of this paragraph.
$p5e = Get-ParaRange("Pour finir le réceptionniste pourra soit valider")
$endRng = $p5e.Duplicate
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1) | Out-Null
$endRng.Collapse(0)
$endRng.Bookmarks.Add("_GoBack") | Out-Null
